# Apply strikethrough formatting to the two "Problem Statement" bullet
# questions about the sprint race correlation and the lap-time correlation.
$d = $word.ActiveDocument

$targets = @(
    "How does the result in the sprint race affect the result in the race?",
    "What is the correlation between the average lap time and fastest lap of a winner and the result of the race?"
)

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text.TrimEnd("`r", "`n", "`x07")
    foreach ($t in $targets) {
        if ($text -eq $t) {
            # Strike through the whole paragraph (including the paragraph
            # mark), matching both the run's rPr and the paragraph's pPr/rPr.
            $p.Range.Font.StrikeThrough = 1
        }
    }
}
